$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 2).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 2).Value = "86.16 (+/-3.52)"
$ws.Cells.Item(5, 2).Value = "79.84 (+/-6.22)"
$ws.Cells.Item(6, 2).Value = "73.94 (+/-8.48)"
$ws.Cells.Item(7, 2).Value = "68.96 (+/-9.61)"
$ws.Cells.Item(8, 2).Value = "65.05 (+/-10.78)"
$ws.Cells.Item(3, 3).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 3).Value = "78.56 (+/-1.50)"
$ws.Cells.Item(5, 3).Value = "66.02 (+/-2.00)"
$ws.Cells.Item(6, 3).Value = "58.64 (+/-2.21)"
$ws.Cells.Item(7, 3).Value = "51.92 (+/-2.30)"
$ws.Cells.Item(8, 3).Value = "44.12 (+/-2.30)"
$ws.Cells.Item(3, 4).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 4).Value = "80.09 (+/-1.27)"
$ws.Cells.Item(5, 4).Value = "66.95 (+/-1.46)"
$ws.Cells.Item(6, 4).Value = "57.13 (+/-1.97)"
$ws.Cells.Item(7, 4).Value = "48.85 (+/-2.31)"
$ws.Cells.Item(8, 4).Value = "40.19 (+/-2.01)"
$ws.Cells.Item(3, 5).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 5).Value = "74.81 (+/-1.62)"
$ws.Cells.Item(5, 5).Value = "61.13 (+/-1.55)"
$ws.Cells.Item(6, 5).Value = "51.01 (+/-1.43)"
$ws.Cells.Item(7, 5).Value = "42.01 (+/-0.96)"
$ws.Cells.Item(8, 5).Value = "35.88 (+/-0.98)"
$ws.Cells.Item(3, 6).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 6).Value = "75.06 (+/-1.02)"
$ws.Cells.Item(5, 6).Value = "61.64 (+/-1.05)"
$ws.Cells.Item(6, 6).Value = "48.25 (+/-2.14)"
$ws.Cells.Item(7, 6).Value = "41.57 (+/-2.00)"
$ws.Cells.Item(8, 6).Value = "35.31 (+/-1.25)"
$ws.Cells.Item(3, 7).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 7).Value = "74.72 (+/-1.84)"
$ws.Cells.Item(5, 7).Value = "60.93 (+/-1.92)"
$ws.Cells.Item(6, 7).Value = "52.33 (+/-2.04)"
$ws.Cells.Item(7, 7).Value = "43.37 (+/-1.71)"
$ws.Cells.Item(8, 7).Value = "37.31 (+/-2.04)"
$ws.Cells.Item(3, 8).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 8).Value = "71.86 (+/-0.91)"
$ws.Cells.Item(5, 8).Value = "57.70 (+/-1.52)"
$ws.Cells.Item(6, 8).Value = "48.52 (+/-1.61)"
$ws.Cells.Item(7, 8).Value = "40.56 (+/-1.94)"
$ws.Cells.Item(8, 8).Value = "32.47 (+/-1.46)"
$ws.Cells.Item(3, 9).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 9).Value = "70.07 (+/-1.29)"
$ws.Cells.Item(5, 9).Value = "54.00 (+/-1.84)"
$ws.Cells.Item(6, 9).Value = "42.13 (+/-1.43)"
$ws.Cells.Item(7, 9).Value = "34.87 (+/-1.15)"
$ws.Cells.Item(8, 9).Value = "26.73 (+/-1.56)"
$ws.Cells.Item(3, 10).Value = "100.00 (+/-0.00)"
$ws.Cells.Item(4, 10).Value = "71.03 (+/-3.09)"
$ws.Cells.Item(5, 10).Value = "57.01 (+/-3.65)"
$ws.Cells.Item(6, 10).Value = "47.41 (+/-3.27)"
$ws.Cells.Item(7, 10).Value = "39.08 (+/-2.51)"
$ws.Cells.Item(8, 10).Value = "31.90 (+/-2.56)"
